$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "36.557.73"
$ws.Cells.Item(2,5).Value = "  +0.50%  "
$ws.Cells.Item(3,4).Value = "1.944.19"
$ws.Cells.Item(3,5).Value = "  -1.03%  "
$ws.Cells.Item(4,5).Value = "  -0.06%  "
$c = $ws.Cells.Item(5,4)
$c.NumberFormat = "@"
$c.Value = "243.56"
$c.Style = "Normal"
$ws.Cells.Item(5,5).Value = "  -0.16%  "
$c = $ws.Cells.Item(6,4)
$c.NumberFormat = "@"
$c.Value = "0.616"
$c.Style = "Normal"
$ws.Cells.Item(6,5).Value = "  +0.62%  "
$ws.Cells.Item(7,5).Value = "  -0.01%  "
$c = $ws.Cells.Item(8,4)
$c.NumberFormat = "@"
$c.Value = "57.56"
$c.Style = "Normal"
$ws.Cells.Item(8,5).Value = "  -1.25%  "
$c = $ws.Cells.Item(9,4)
$c.NumberFormat = "@"
$c.Value = "0.362"
$c.Style = "Normal"
$ws.Cells.Item(9,5).Value = "  -2.39%  "
$c = $ws.Cells.Item(10,4)
$c.NumberFormat = "@"
$c.Value = "0.0847"
$c.Style = "Normal"
$ws.Cells.Item(10,5).Value = "  -0.11%  "
$ws.Cells.Item(11,5).Value = "  -0.20%  "
$ws.Cells.Item(12,4).Value = "2.228.36"
$ws.Cells.Item(13,2).Value = "Polygon"
$ws.Cells.Item(13,3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$c = $ws.Cells.Item(13,4)
$c.NumberFormat = "@"
$c.Value = "0.813"
$c.Style = "Normal"
$ws.Cells.Item(13,5).Value = "  -2.75%  "
$ws.Cells.Item(14,2).Value = "Avalanche"
$ws.Cells.Item(14,3).Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$c = $ws.Cells.Item(14,4)
$c.NumberFormat = "@"
$c.Value = "21.30"
$c.Style = "Normal"
$ws.Cells.Item(14,5).Value = "  -4.69%  "
$c = $ws.Cells.Item(15,4)
$c.NumberFormat = "@"
$c.Value = "13.47"
$c.Style = "Normal"
$ws.Cells.Item(15,5).Value = "  -0.38%  "
$ws.Cells.Item(16,5).Value = "  -3.58%  "
$ws.Cells.Item(17,4).Value = "1.944.44"
$ws.Cells.Item(17,5).Value = "  -1.76%  "
$ws.Cells.Item(18,4).Value = "36.445.84"
$ws.Cells.Item(18,5).Value = "  +0.45%  "
$c = $ws.Cells.Item(19,4)
$c.NumberFormat = "@"
$c.Value = "69.42"
$c.Style = "Normal"
$ws.Cells.Item(19,5).Value = "  -2.83%  "
$ws.Cells.Item(20,4).Value = "0.0₃0865"
$ws.Cells.Item(20,5).Value = "  -2.28%  "
$c = $ws.Cells.Item(21,4)
$c.NumberFormat = "@"
$c.Value = "228.62"
$c.Style = "Normal"
$ws.Cells.Item(21,5).Value = "  -0.51%  "
$c = $ws.Cells.Item(22,4)
$c.NumberFormat = "@"
$c.Value = "5.00"
$c.Style = "Normal"
$ws.Cells.Item(22,5).Value = "  -3.01%  "
$ws.Cells.Item(23,5).Value = "  +0.01%  "
$c = $ws.Cells.Item(24,4)
$c.NumberFormat = "@"
$c.Value = "2.36"
$c.Style = "Normal"
$ws.Cells.Item(24,5).Value = "  -6.32%  "
$ws.Cells.Item(25,5).Value = "  +1.73%  "
$c = $ws.Cells.Item(26,4)
$c.NumberFormat = "@"
$c.Value = "9.22"
$c.Style = "Normal"
$ws.Cells.Item(26,5).Value = "  -3.27%  "
$c = $ws.Cells.Item(27,4)
$c.NumberFormat = "@"
$c.Value = "161.73"
$c.Style = "Normal"
$ws.Cells.Item(27,5).Value = "  -2.45%  "
$ws.Cells.Item(28,5).Value = "  +8.40%  "
$c = $ws.Cells.Item(29,4)
$c.NumberFormat = "@"
$c.Value = "19.20"
$c.Style = "Normal"
$ws.Cells.Item(29,5).Value = "  -3.84%  "
$ws.Cells.Item(30,5).Value = "  -0.78%  "
$ws.Cells.Item(31,5).Value = "  -4.73%  "
$c = $ws.Cells.Item(32,4)
$c.NumberFormat = "@"
$c.Value = "4.60"
$c.Style = "Normal"
$ws.Cells.Item(32,5).Value = "  -3.47%  "
$ws.Cells.Item(33,5).Value = "  -3.91%  "
$c = $ws.Cells.Item(34,4)
$c.NumberFormat = "@"
$c.Value = "4.20"
$c.Style = "Normal"
$ws.Cells.Item(34,5).Value = "  -3.26%  "
$c = $ws.Cells.Item(35,4)
$c.NumberFormat = "@"
$c.Value = "6.24"
$c.Style = "Normal"
$ws.Cells.Item(35,5).Value = "  +4.63%  "
$ws.Cells.Item(36,5).Value = "  -0.16%  "
$ws.Cells.Item(37,5).Value = "  -1.18%  "
$ws.Cells.Item(38,5).Value = "  +0.88%  "
$c = $ws.Cells.Item(39,4)
$c.NumberFormat = "@"
$c.Value = "3.17"
$c.Style = "Normal"
$ws.Cells.Item(39,5).Value = "  +8.60%  "
$c = $ws.Cells.Item(40,4)
$c.NumberFormat = "@"
$c.Value = "0.0983"
$c.Style = "Normal"
$ws.Cells.Item(40,5).Value = "  +1.98%  "
$ws.Cells.Item(41,5).Value = "  +0.42%  "
$ws.Cells.Item(42,2).Value = "TrustWalletToken"
$ws.Cells.Item(42,3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$c = $ws.Cells.Item(42,4)
$c.NumberFormat = "@"
$c.Value = "1.15"
$c.Style = "Normal"
$ws.Cells.Item(42,5).Value = "  -2.61%  "
$ws.Cells.Item(43,2).Value = "VeChain"
$ws.Cells.Item(43,3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Cells.Item(43,4)
$c.NumberFormat = "@"
$c.Value = "0.0209"
$c.Style = "Normal"
$ws.Cells.Item(43,5).Value = "  -0.96%  "
$c = $ws.Cells.Item(44,4)
$c.NumberFormat = "@"
$c.Value = "16.03"
$c.Style = "Normal"
$ws.Cells.Item(44,5).Value = "  +1.25%  "
$ws.Cells.Item(45,4).Value = "1.345.03"
$ws.Cells.Item(45,5).Value = "  -0.38%  "
$ws.Cells.Item(46,5).Value = "  -2.62%  "
$c = $ws.Cells.Item(47,4)
$c.NumberFormat = "@"
$c.Value = "86.75"
$c.Style = "Normal"
$ws.Cells.Item(47,5).Value = "  -2.43%  "
$c = $ws.Cells.Item(48,4)
$c.NumberFormat = "@"
$c.Value = "7.19"
$c.Style = "Normal"
$ws.Cells.Item(48,5).Value = "  -1.39%  "
$ws.Cells.Item(49,5).Value = "  -0.09%  "
$ws.Cells.Item(50,4).Value = "2.120.43"
$ws.Cells.Item(50,5).Value = "  -1.02%  "
$c = $ws.Cells.Item(51,4)
$c.NumberFormat = "@"
$c.Value = "43.35"
$c.Style = "Normal"
$ws.Cells.Item(51,5).Value = "  -4.12%  "
